$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (forecast values) and column C (hour values) for the
# affected rows. Rows 8, 9, and 16 are left untouched.

$ws.Range("B2").Value = 318.5
$ws.Range("C2").Value = 38

$ws.Range("B3").Value = 206
$ws.Range("C3").Value = 38

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 38

$ws.Range("B5").Value = 106
$ws.Range("C5").Value = 38

$ws.Range("B6").Value = 38
$ws.Range("C6").Value = 38

$ws.Range("B7").Value = 82
$ws.Range("C7").Value = 38

$ws.Range("B10").Value = 379
$ws.Range("C10").Value = 38

$ws.Range("B11").Value = 208
$ws.Range("C11").Value = 38

$ws.Range("B12").Value = 395.2
$ws.Range("C12").Value = 38

$ws.Range("B13").Value = 302
$ws.Range("C13").Value = 38

$ws.Range("B14").Value = 495
$ws.Range("C14").Value = 38

$ws.Range("B15").Value = 127
$ws.Range("C15").Value = 38

$ws.Range("B17").Value = 77
$ws.Range("C17").Value = 38

$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 38

$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 38

$ws.Range("B20").Value = 38
$ws.Range("C20").Value = 38
